$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Fecha(D), Volumen(M), Precio minimo(N), Precio maximo(O), Precio promedio ponderado(P),
# Unidad de comercializacion(Q), Origen(R), Precio $/Kg(S), Kg/unidad(T)
$data = @{
    2 = @{ D = 44585; M = 50; N = 22500; O = 22500; P = 22500; Q = '$/caja 15 kilos empedrada'; R = 'Provincia de Limarí'; S = 1500; T = 15 }
    3 = @{ D = 44239; M = 70; N = 15000; O = 15000; P = 15000; Q = '$/caja 15 kilos granel'; R = 'Provincia de Limarí'; S = 1000; T = 15 }
    4 = @{ D = 44259; M = 80; N = 12000; O = 12000; P = 12000; Q = '$/caja 15 kilos empedrada'; R = 'Provincia de Limarí'; S = 800; T = 15 }
    5 = @{ D = 44270; M = 85; N = 12000; O = 12000; P = 12000; Q = '$/caja 14 kilos granel'; R = 'Provincia del Elquí'; S = 857; T = 14 }
    6 = @{ D = 44278; M = 45; N = 13000; O = 13000; P = 13000; Q = '$/caja 14 kilos empedrada'; R = 'Provincia del Elquí'; S = 929; T = 14 }
    7 = @{ D = 44314; M = 56; N = 14000; O = 14000; P = 14000; Q = '$/caja 14 kilos granel'; R = 'Provincia de Limarí'; S = 1000; T = 14 }
    8 = @{ D = 44260; M = 56; N = 13000; O = 13000; P = 13000; Q = '$/caja 14 kilos empedrada'; R = 'Provincia del Elquí'; S = 929; T = 14 }
    9 = @{ D = 44245; M = 50; N = 15000; O = 15000; P = 15000; Q = '$/caja 15 kilos granel'; R = 'Provincia de Limarí'; S = 1000; T = 15 }
    10 = @{ D = 44323; M = 60; N = 14000; O = 14000; P = 14000; Q = '$/caja 14 kilos granel'; R = 'Provincia de Limarí'; S = 1000; T = 14 }
    11 = @{ D = 44592; M = 54; N = 20000; O = 20000; P = 20000; Q = '$/caja 15 kilos empedrada'; R = 'Provincia de Limarí'; S = 1333; T = 15 }
    12 = @{ D = 44316; M = 48; N = 14000; O = 14000; P = 14000; Q = '$/caja 14 kilos granel'; R = 'Provincia de Limarí'; S = 1000; T = 14 }
    13 = @{ D = 44322; M = 50; N = 14000; O = 14000; P = 14000; Q = '$/caja 14 kilos granel'; R = 'Provincia de Limarí'; S = 1000; T = 14 }
    14 = @{ D = 44588; M = 85; N = 19000; O = 20000; P = 19529; Q = '$/caja 14 kilos granel'; R = 'Provincia de Limarí'; S = 1395; T = 14 }
    15 = @{ D = 44320; M = 45; N = 14000; O = 14000; P = 14000; Q = '$/caja 14 kilos granel'; R = 'Provincia de Limarí'; S = 1000; T = 14 }
    16 = @{ D = 44313; M = 36; N = 14000; O = 14000; P = 14000; Q = '$/caja 14 kilos granel'; R = 'Provincia de Limarí'; S = 1000; T = 14 }
    17 = @{ D = 44238; M = 60; N = 15000; O = 15000; P = 15000; Q = '$/caja 15 kilos granel'; R = 'Provincia de Limarí'; S = 1000; T = 15 }
    18 = @{ D = 44242; M = 45; N = 12000; O = 12000; P = 12000; Q = '$/caja 15 kilos granel'; R = 'Provincia de Limarí'; S = 800; T = 15 }
    19 = @{ D = 44252; M = 60; N = 14000; O = 14000; P = 14000; Q = '$/caja 14 kilos empedrada'; R = 'Provincia de Limarí'; S = 1000; T = 14 }
    20 = @{ D = 44271; M = 50; N = 12000; O = 12000; P = 12000; Q = '$/caja 14 kilos granel'; R = 'Provincia del Elquí'; S = 857; T = 14 }
    21 = @{ D = 44315; M = 65; N = 14000; O = 14000; P = 14000; Q = '$/caja 14 kilos granel'; R = 'Provincia de Limarí'; S = 1000; T = 14 }
    22 = @{ D = 44312; M = 68; N = 14000; O = 14000; P = 14000; Q = '$/caja 14 kilos granel'; R = 'Provincia de Limarí'; S = 1000; T = 14 }
}

foreach ($row in $data.Keys) {
    $v = $data[$row]
    $ws.Cells.Item($row, 4).Value = $v.D
    $ws.Cells.Item($row, 13).Value = $v.M
    $ws.Cells.Item($row, 14).Value = $v.N
    $ws.Cells.Item($row, 15).Value = $v.O
    $ws.Cells.Item($row, 16).Value = $v.P
    $ws.Cells.Item($row, 17).Value = $v.Q
    $ws.Cells.Item($row, 18).Value = $v.R
    $ws.Cells.Item($row, 19).Value = $v.S
    $ws.Cells.Item($row, 20).Value = $v.T
}
